$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: human-readable, capitalized labels
$ws.Range("A1").Value = "Edad"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Situación profesional código"
$ws.Range("D1").Value = "Aragón"
$ws.Range("E1").Value = "Situación profesional"
$ws.Range("F1").Value = "Sexo"

# Row 2: measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:edad"
$ws.Range("B2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:situacion-profesional"
$ws.Range("F2").Value = "iaest-measure:sexo"

# Row 3: medida/dim/null markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Row 4: xsd type / codelist identifiers
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "xsd:string"
